# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (totals) sheet,
#    populated with the same fund-holding layout used by the other quarter
#    sheets (header row + two currency variants of the same fund).
# 2. Insert a new leading data row into "总计" summarizing the 2022-Q1
#    quarter, shifting the existing rows down and renumbering the index
#    column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Borrow the header / index-column formatting from the neighbouring
# "2021-Q4" sheet so the new sheet matches the look of the rest of the
# workbook (bold, centered, top-aligned, thin border).
$reference = $wb.Worksheets.Item("2021-Q4")
$reference.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$reference.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - fund 000369
$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "000369"
$q1.Range("B2").Style = "Normal"
$q1.Range("C2").Value = "广发全球医疗保健(QDII) - 人民币"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "2.46"
$q1.Range("D2").Style = "Normal"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "81.85"
$q1.Range("E2").Style = "Normal"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "3.72"
$q1.Range("F2").Style = "Normal"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0915"
$q1.Range("G2").Style = "Normal"
$q1.Range("H2").Value = 3

# Row 3 - fund 000370
$q1.Range("A3").Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "000370"
$q1.Range("B3").Style = "Normal"
$q1.Range("C3").Value = "广发全球医疗保健(QDII) - 美元"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "2.46"
$q1.Range("D3").Style = "Normal"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "81.85"
$q1.Range("E3").Style = "Normal"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "3.72"
$q1.Range("F3").Style = "Normal"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.0915"
$q1.Range("G3").Style = "Normal"
$q1.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2) Prepend a 2022-Q1 summary row into "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

# Copy the index-column formatting from the row below (now row 3) so the
# new A2 cell matches the look of the rest of the index column.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.18
$total.Range("B2:D2").Style = "Normal"

# Renumber the remaining index column (previously 0..4, now 1..5)
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Restore the original active sheet/selection (adding sheets moves it).
$original = $wb.Worksheets.Item("2020-Q4")
$original.Activate() | Out-Null
$original.Range("A1").Select() | Out-Null
